$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","G","H","I","J","K")

$data = @(
  @(-19.54155041195853, 4.321926298890479, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853),
  @(-19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, 2.288722224122679, -19.54155041195853, -19.54155041195853),
  @(-19.54155041195853, -19.54155041195853, 2.777579496263461, -19.54155041195853, 2.436006579013743, -19.54155041195853, 1.790873554221308, -19.54155041195853, 2.354690700022308, -19.54155041195853),
  @(-19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, 2.10167208209437, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853),
  @(-19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853),
  @(3.112099630502672, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853),
  @(-19.54155041195853, -19.54155041195853, -19.54155041195853, 3.055365672919507, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853),
  @(3.505074878883625, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853),
  @(-19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, 1.502692271333988, -19.54155041195853, 2.238576403520431),
  @(-19.54155041195853, -19.54155041195853, -19.54155041195853, 1.707751149530495, -19.54155041195853, 2.577608881583151, -19.54155041195853, -19.54155041195853, -19.54155041195853, 1.385423800728065),
  @(-19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853),
  @(-19.54155041195853, -19.54155041195853, -19.54155041195853, 1.464935244263216, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, 2.283240946178962, 1.618474243415421),
  @(-19.54155041195853, -19.54155041195853, 1.755532096113401, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, 2.0814277002101),
  @(-19.54155041195853, -19.54155041195853, -1.775765766098652, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853),
  @(-19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, 2.322639133962368, -19.54155041195853),
  @(-19.54155041195853, -19.54155041195853, -0.9303202423976358, -19.54155041195853, -19.54155041195853, -19.54155041195853, 0.4814688618291562, 0.9007583401555037, 1.288219543120929, -19.54155041195853),
  @(-19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, -19.54155041195853, 0.357247275646322, 0.9015884113292372, 1.363284834649255, -19.54155041195853),
  @(-19.54155041195853, -19.54155041195853, 1.873589413368655, -19.54155041195853, -19.54155041195853, -19.54155041195853, 1.955771779165483, 2.17137829058335, -19.54155041195853, -19.54155041195853),
  @(-19.54155041195853, -19.54155041195853, 2.401972079644047, -19.54155041195853, 3.866765462815013, -19.54155041195853, 2.20817378366108, 2.014475177030621, -19.54155041195853, 2.423691703540175),
  @(-19.54155041195853, -19.54155041195853, -19.54155041195853, 2.500795973981129, -19.54155041195853, 3.283687320783918, 2.422895166227352, -19.54155041195853, -19.54155041195853, -19.54155041195853)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $addr = $cols[$j] + $rowNum
        $ws.Range($addr).Value = $rowData[$j]
    }
}
